$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet2 ("Sheet1") - build the WGS / Django-model-field overview table
# ---------------------------------------------------------------------------

# Column A: field names (these are all pre-existing shared strings already
# used as column headers on the Antibiotics sheet, so no new sharedString
# entries are created here).
$ws2.Range("A1").Value = "Show"
$ws2.Range("A2").Value = "Retest"
$ws2.Range("A3").Value = "Disk_Abx"
$ws2.Range("A4").Value = "Test_Method"
$ws2.Range("A5").Value = "Abx_Code"
$ws2.Range("A6").Value = "Whonet_Abx"
$ws2.Range("A7").Value = "Antibiotic"
$ws2.Range("A8").Value = "Guidelines"
$ws2.Range("A9").Value = "Potency"
$ws2.Range("A10").Value = "Class"
$ws2.Range("A11").Value = "Subclass"
$ws2.Range("A12").Value = "Date_Modified"

# Column B / C: model field type + field kwargs (order matters: this is the
# order in which brand-new shared strings get appended to the workbook).
$ws2.Range("B1").Value = "BooleanField"
$ws2.Range("B12").Value = "DateField"
$ws2.Range("C1").Value = "default=True"
$ws2.Range("C12").Value = "auto_now_add=True"
$ws2.Range("C4").Value = "max_length=100, blank=True, default=''''"

# Column E: literal (typed / pasted) description text for each row.
$ws2.Range("E1").Value = "Show=models.BooleanField(default=True)"
$ws2.Range("E2").Value = "Retest=models.BooleanField(default=True)"
$ws2.Range("E3").Value = "Disk_Abx=models.BooleanField(default=True)"
$ws2.Range("E12").Value = "Date_Modified=models.DateField(auto_now_add=True)"
$ws2.Range("E4").Value = 'Test_Method=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E5").Value = 'Abx_Code=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E6").Value = 'Whonet_Abx=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E7").Value = 'Antibiotic=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E8").Value = 'Guidelines=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E9").Value = 'Potency=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E10").Value = 'Class=models.Charfield(max_length=100, blank=True, default="")'
$ws2.Range("E11").Value = 'Subclass=models.Charfield(max_length=100, blank=True, default="")'

# Column B rows 4-11 all share the "CharField" type text - last new string.
$ws2.Range("B4:B11").Value = "CharField"

# Fill in C5:C11 to match C4 (same kwargs text, reuses the shared string).
$ws2.Range("C5:C11").Value = "max_length=100, blank=True, default=''''"
$ws2.Range("B2:B3").Value = "BooleanField"
$ws2.Range("C2:C3").Value = "default=True"

# Column D: CONCATENATE formula building the Django field declaration.
# D1 stands alone; D2:D12 form a shared-formula group.
$ws2.Range("D1").Formula = "=CONCATENATE(A1,""="",""models."",B1,""("",C1,"")"","""")"
$ws2.Range("D2:D12").Formula = "=CONCATENATE(A2,""="",""models."",B2,""("",C2,"")"","""")"

# Header-row-style formatting pulled from the Antibiotics sheet but
# left-aligned instead of centered.
$ws1.Range("A1").Copy()
$ws2.Range("A1:A4").PasteSpecial(-4122)
$ws2.Range("A1:A4").HorizontalAlignment = -4131

$ws1.Range("E1").Copy()
$ws2.Range("A5:A12").PasteSpecial(-4122)
$ws2.Range("A5:A12").HorizontalAlignment = -4131

# Column widths (best-fit approximations).
$ws2.Columns.Item(1).ColumnWidth = 13.6
$ws2.Columns.Item(2).ColumnWidth = 11.8
$ws2.Columns.Item(3).ColumnWidth = 35.6
$ws2.Range("D1:E1").ColumnWidth = 67.6

# ---------------------------------------------------------------------------
# Sheet selection / active-tab bookkeeping: the author ended up on Sheet1
# with A1:L1 selected on the Antibiotics sheet.
# ---------------------------------------------------------------------------
$ws1.Range("A1:L1").Select()
$ws2.Range("E22").Select()
$ws2.Activate()
